# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns
# to match the latest scrape. D-column writes force the NumberFormat to
# text ("@") before assigning so Excel does not auto-convert numeric-
# looking price strings (e.g. "1.00" -> 1, "600.35" -> 600.35000000000002)
# into Number cells, then reapply the "Normal" style so no stray style
# index is introduced (the source cells carry no explicit style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.668.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.00%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.665.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.614"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.124"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.38%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.45%  "

$ws.Range("E11").Value = "  -0.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.154"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000196"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.95%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.147.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.466.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.666.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000110"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.78%  "

$ws.Range("E26").Value = "  -5.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.167"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("E28").Value = "  -3.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "542.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.59%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.13%  "

$ws.Range("E32").Value = "  -3.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.46%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.421"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.72%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.57%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "159.26"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.02%  "

$ws.Range("E40").Value = "  -2.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.82%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "164.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0612"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.643"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.50%  "

$ws.Range("E49").Value = "  -2.04%  "

$ws.Range("E50").Value = "  +3.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.08%  "
